# Applies the "handles float input without breaking stuff" marksheet update:
#  - Summary block (rows 10-12): real Right/Wrong/NotAttempt/Max counts,
#    numeric "-1" marking penalty (was stored as text), and a "86/112"
#    total instead of the placeholder "Absent".
#  - Per-question Student-Ans columns (A/D) are now filled in with the
#    option the student actually picked (highlighted correct/incorrect),
#    and the extra third G/H "Student Ans / Correct Ans" block is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block ----------------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "86/112"

# ---- Drop the third (G/H) Student Ans / Correct Ans block -------------
$ws.Range("G15:H40").Clear()

# ---- Fill in the Student Ans values for column A (and D where it still
#      exists, rows 16-18), styling each as correct or incorrect --------
function Set-StudentAns($cellRef, $styleName, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Style = $styleName
    $rng.Value = $text
}

Set-StudentAns "A16" "correctStyle" "Option A"
Set-StudentAns "D16" "correctStyle" "Option A"

Set-StudentAns "A17" "correctStyle" "Option D"
Set-StudentAns "D17" "correctStyle" "Option C"

Set-StudentAns "A18" "correctStyle" "Option B"
Set-StudentAns "D18" "correctStyle" "Option D"

Set-StudentAns "A19" "correctStyle" "Option C"
$ws.Range("D19:E19").Clear()

$ws.Range("D20:E20").Clear()

Set-StudentAns "A21" "correctStyle" "Option C"
$ws.Range("D21:E21").Clear()

Set-StudentAns "A22" "correctStyle" "Option D"
$ws.Range("D22:E22").Clear()

Set-StudentAns "A23" "correctStyle" "Option D"
$ws.Range("D23:E23").Clear()

$ws.Range("D24:E24").Clear()

Set-StudentAns "A25" "correctStyle" "Option A"
$ws.Range("D25:E25").Clear()

Set-StudentAns "A26" "correctStyle" "Option C"
$ws.Range("D26:E26").Clear()

Set-StudentAns "A27" "correctStyle" "Option A"
$ws.Range("D27:E27").Clear()

Set-StudentAns "A28" "correctStyle" "Option D"
$ws.Range("D28:E28").Clear()

Set-StudentAns "A29" "correctStyle" "Option D"
$ws.Range("D29:E29").Clear()

Set-StudentAns "A30" "correctStyle" "Option B"
$ws.Range("D30:E30").Clear()

Set-StudentAns "A31" "correctStyle" "Option D"
$ws.Range("D31:E31").Clear()

Set-StudentAns "A32" "correctStyle" "Option C"
$ws.Range("D32:E32").Clear()

Set-StudentAns "A33" "correctStyle" "Option D"
$ws.Range("D33:E33").Clear()

Set-StudentAns "A34" "incorrectStyle" "Option A"
$ws.Range("D34:E34").Clear()

Set-StudentAns "A35" "correctStyle" "Option D"
$ws.Range("D35:E35").Clear()

Set-StudentAns "A36" "incorrectStyle" "Option D"
$ws.Range("D36:E36").Clear()

$ws.Range("D37:E37").Clear()

Set-StudentAns "A38" "correctStyle" "Option A"
$ws.Range("D38:E38").Clear()

Set-StudentAns "A39" "correctStyle" "Option D"
$ws.Range("D39:E39").Clear()

$ws.Range("D40:E40").Clear()

Write-Output "edit applied"
